$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.850.62'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.736.73'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '233.90'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.5176'
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("D8").Value = '0.2766'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").Value = '39.28'
$ws.Range("E9").Value = '  -2.70%  '
$ws.Range("D10").Value = '0.06122'
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '1.736.84'
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = '0.07048'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '15.25'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").Value = '0.6391'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").Value = '4.509'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = '76.86'
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").Value = '0.9990'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").Value = '25.831.62'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").Value = '11.48'
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").Value = '0.000006640'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").Value = '1.957.97'
$ws.Range("E22").Value = '  -1.97%  '
$ws.Range("D23").Value = '4.137'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").Value = '8.748'
$ws.Range("E24").Value = '  +5.73%  '
$ws.Range("D25").Value = '5.136'
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").Value = '139.79'
$ws.Range("E26").Value = '  +2.62%  '
$ws.Range("D27").Value = '1.514'
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("D28").Value = '15.01'
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").Value = '1.787'
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").Value = '101.90'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '0.08301'
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").Value = '3.690'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '3.450'
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").Value = '0.04507'
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").Value = '2.612'
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").Value = '0.9772'
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").Value = '0.6129'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").Value = '2.661'
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").Value = '0.01581'
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("D40").Value = '1.943'
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").Value = '0.9986'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '100.61'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").Value = '0.3836'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '0.7248'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("D45").Value = '4.979'
$ws.Range("E45").Value = '  +1.78%  '
$ws.Range("D46").Value = '0.05377'
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("D47").Value = '0.1125'
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("D48").Value = '6.249'
$ws.Range("E48").Value = '  +4.56%  '
$ws.Range("D49").Value = '52.98'
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").Value = '29.98'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '7.577'
$ws.Range("E51").Value = '  +1.94%  '
